$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift existing BOM rows 6-15 down to 7-16 to make room for the new fuse
# row, carrying formatting along (avoids Rows.Insert(), which would pull in
# extra unused cell styles).
for ($r = 15; $r -ge 6; $r--) {
    $src = $ws.Range("A" + $r + ":D" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $src.Copy()
    $dst.PasteSpecial(-4163)
}

# Clear the old contents left behind in row 6 before writing the new part
$ws.Range("A6:D6").ClearContents()

# New BOM entry: 1A resettable fuse (F1) added to the 12V power rail
$ws.Cells.Item(6, 2).Value = "F1"
$ws.Cells.Item(6, 3).Value = 1206
$ws.Cells.Item(6, 4).Value = "C70081"

$ws.Range("D21").Select()
